# Update gh-pages data values (generated output refresh).
# Updates "想去人数" (interest count) and "最低票价" (min price) figures
# across the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 238
$ws.Range("F5").Value = 13952
$ws.Range("G5").Value = 70
$ws.Range("F8").Value = 221
$ws.Range("F14").Value = 542
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = 14021
$ws.Range("F21").Value = 14994
$ws.Range("F23").Value = 8305
$ws.Range("F24").Value = 280
$ws.Range("F26").Value = 30
$ws.Range("F27").Value = 157
$ws.Range("F30").Value = 4
$ws.Range("F39").Value = 10
$ws.Range("F42").Value = 221
$ws.Range("F44").Value = 99
$ws.Range("F45").Value = 5112

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 51

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 238
$ws.Range("F5").Value = 13952
$ws.Range("G5").Value = 70
$ws.Range("F8").Value = 221
$ws.Range("F14").Value = 542
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = 14021
$ws.Range("F21").Value = 14994
$ws.Range("F23").Value = 8305
$ws.Range("F24").Value = 280
$ws.Range("F26").Value = 30
$ws.Range("F27").Value = 157
$ws.Range("F30").Value = 4
$ws.Range("F38").Value = 51
$ws.Range("F41").Value = 10
$ws.Range("F44").Value = 221
$ws.Range("F46").Value = 99
$ws.Range("F47").Value = 5112
